# Add two new column-name / category rows to the metadata mapping sheet
# ("Added columns to clean"): cell_authors_annotation -> annotation_authors
# and subtype -> disease_subtype, appended right after the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "cell_authors_annotation"
$ws.Range("B13").Value = "annotation_authors"

$ws.Range("A14").Value = "subtype"
$ws.Range("B14").Value = "disease_subtype"
